$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lama2"
$ws.Range("C2").Value = "Itga7"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.471482
$ws.Range("H2").Value = 4.414446
$ws.Range("I2").Value = 0.004946458467382327
$ws.Range("J2").Value = 0.004946458467382326
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.464265666666667
$ws.Range("N2").Value = 4.392797
$ws.Range("O2").Value = 0.02620474750556022
$ws.Range("P2").Value = 0.02620474750556022
$ws.Range("Q2").Value = 2.154640571718
$ws.Range("R2").Value = 19.391765145462
$ws.Range("S2").Value = 0.0001296206951844943
$ws.Range("T2").Value = 0.0001296206951844942

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lama2"
$ws.Range("C3").Value = "Itga7"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.471482
$ws.Range("H3").Value = 4.414446
$ws.Range("I3").Value = 0.004946458467382327
$ws.Range("J3").Value = 0.004946458467382326
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8090393333333333
$ws.Range("N3").Value = 2.427118
$ws.Range("O3").Value = 0.01447870556190061
$ws.Range("P3").Value = 0.01447870556190061
$ws.Range("Q3").Value = 1.190486816292
$ws.Range("R3").Value = 10.714381346628
$ws.Range("S3").Value = 0.00007161831572339887
$ws.Range("T3").Value = 0.00007161831572339885

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lama2"
$ws.Range("C4").Value = "Itga7"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.471482
$ws.Range("H4").Value = 4.414446
$ws.Range("I4").Value = 0.004946458467382327
$ws.Range("J4").Value = 0.004946458467382326
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 53.60457233333333
$ws.Range("N4").Value = 160.813717
$ws.Range("O4").Value = 0.9593165469325392
$ws.Range("P4").Value = 0.9593165469325391
$ws.Range("Q4").Value = 78.87816330619799
$ws.Range("R4").Value = 709.903469755782
$ws.Range("S4").Value = 0.004745219456474434
$ws.Range("T4").Value = 0.004745219456474433

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lama2"
$ws.Range("C5").Value = "Itga7"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 264.7713316666666
$ws.Range("H5").Value = 794.313995
$ws.Range("I5").Value = 0.8900417371348598
$ws.Range("J5").Value = 0.8900417371348596
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.464265666666667
$ws.Range("N5").Value = 4.392797
$ws.Range("O5").Value = 0.02620474750556022
$ws.Range("P5").Value = 0.02620474750556022
$ws.Range("Q5").Value = 387.6955704771128
$ws.Range("R5").Value = 3489.260134294015
$ws.Range("S5").Value = 0.02332331899102921
$ws.Range("T5").Value = 0.0233233189910292

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lama2"
$ws.Range("C6").Value = "Itga7"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 264.7713316666666
$ws.Range("H6").Value = 794.313995
$ws.Range("I6").Value = 0.8900417371348598
$ws.Range("J6").Value = 0.8900417371348596
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8090393333333333
$ws.Range("N6").Value = 2.427118
$ws.Range("O6").Value = 0.01447870556190061
$ws.Range("P6").Value = 0.01447870556190061
$ws.Range("Q6").Value = 214.2104216573789
$ws.Range("R6").Value = 1927.89379491641
$ws.Range("S6").Value = 0.01288665224977818
$ws.Range("T6").Value = 0.01288665224977817

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Lama2"
$ws.Range("C7").Value = "Itga7"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 264.7713316666666
$ws.Range("H7").Value = 794.313995
$ws.Range("I7").Value = 0.8900417371348598
$ws.Range("J7").Value = 0.8900417371348596
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 53.60457233333333
$ws.Range("N7").Value = 160.813717
$ws.Range("O7").Value = 0.9593165469325392
$ws.Range("P7").Value = 0.9593165469325391
$ws.Range("Q7").Value = 14192.95400011882
$ws.Range("R7").Value = 127736.5860010694
$ws.Range("S7").Value = 0.8538317658940524
$ws.Range("T7").Value = 0.8538317658940522

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Lama2"
$ws.Range("C8").Value = "Itga7"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 31.239114
$ws.Range("H8").Value = 93.717342
$ws.Range("I8").Value = 0.105011804397758
$ws.Range("J8").Value = 0.105011804397758
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.464265666666667
$ws.Range("N8").Value = 4.392797
$ws.Range("O8").Value = 0.02620474750556022
$ws.Range("P8").Value = 0.02620474750556022
$ws.Range("Q8").Value = 45.742362087286
$ws.Range("R8").Value = 411.681258785574
$ws.Range("S8").Value = 0.002751807819346528
$ws.Range("T8").Value = 0.002751807819346527

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Lama2"
$ws.Range("C9").Value = "Itga7"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 31.239114
$ws.Range("H9").Value = 93.717342
$ws.Range("I9").Value = 0.105011804397758
$ws.Range("J9").Value = 0.105011804397758
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8090393333333333
$ws.Range("N9").Value = 2.427118
$ws.Range("O9").Value = 0.01447870556190061
$ws.Range("P9").Value = 0.01447870556190061
$ws.Range("Q9").Value = 25.273671964484
$ws.Range("R9").Value = 227.463047680356
$ws.Range("S9").Value = 0.001520434996399039
$ws.Range("T9").Value = 0.001520434996399038

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Lama2"
$ws.Range("C10").Value = "Itga7"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 31.239114
$ws.Range("H10").Value = 93.717342
$ws.Range("I10").Value = 0.105011804397758
$ws.Range("J10").Value = 0.105011804397758
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 53.60457233333333
$ws.Range("N10").Value = 160.813717
$ws.Range("O10").Value = 0.9593165469325392
$ws.Range("P10").Value = 0.9593165469325391
$ws.Range("Q10").Value = 1674.559346042246
$ws.Range("R10").Value = 15071.03411438021
$ws.Range("S10").Value = 0.1007395615820125
$ws.Range("T10").Value = 0.1007395615820124
